{"js": "// Update the date line and the 25 multiplication-fact answers in the\n// practice table. Each old value is unique in the document, so a\n// search-and-replace keyed on the exact old text is unambiguous and\n// preserves the existing run formatting (fonts/size) on every cell.\nconst replacements = [\n  [\"2025-06-24 Tuesday\", \"2025-06-25 Wednesday\"],\n  [\"99\u00d767=6633\", \"25\u00d750=1250\"],\n  [\"47\u00d792=4324\", \"94\u00d739=3666\"],\n  [\"82\u00d757=4674\", \"39\u00d753=2067\"],\n  [\"69\u00d763=4347\", \"72\u00d781=5832\"],\n  [\"68\u00d796=6528\", \"12\u00d752=624\"],\n  [\"61\u00d799=6039\", \"45\u00d741=1845\"],\n  [\"46\u00d742=1932\", \"64\u00d748=3072\"],\n  [\"68\u00d714=952\", \"22\u00d769=1518\"],\n  [\"50\u00d727=1350\", \"76\u00d743=3268\"],\n  [\"95\u00d784=7980\", \"25\u00d723=575\"],\n  [\"69\u00d779=5451\", \"95\u00d716=1520\"],\n  [\"86\u00d787=7482\", \"49\u00d722=1078\"],\n  [\"92\u00d790=8280\", \"56\u00d740=2240\"],\n  [\"74\u00d740=2960\", \"41\u00d796=3936\"],\n  [\"35\u00d771=2485\", \"36\u00d717=612\"],\n  [\"57\u00d742=2394\", \"30\u00d745=1350\"],\n  [\"35\u00d722=770\", \"95\u00d720=1900\"],\n  [\"38\u00d724=912\", \"63\u00d767=4221\"],\n  [\"69\u00d790=6210\", \"17\u00d778=1326\"],\n  [\"93\u00d740=3720\", \"44\u00d727=1188\"],\n  [\"45\u00d768=3060\", \"35\u00d756=1960\"],\n  [\"43\u00d745=1935\", \"58\u00d741=2378\"],\n  [\"39\u00d797=3783\", \"61\u00d789=5429\"],\n  [\"23\u00d724=552\", \"77\u00d790=6930\"],\n  [\"49\u00d742=2058\", \"17\u00d784=1428\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 multiplication-fact answers in the\n# practice table. Each old value is unique in the document, so a\n# Find/Replace keyed on the exact old text is unambiguous and preserves\n# the existing run formatting (fonts/size) on every cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-24 Tuesday\", \"2025-06-25 Wednesday\"),\n    @(\"99\u00d767=6633\", \"25\u00d750=1250\"),\n    @(\"47\u00d792=4324\", \"94\u00d739=3666\"),\n    @(\"82\u00d757=4674\", \"39\u00d753=2067\"),\n    @(\"69\u00d763=4347\", \"72\u00d781=5832\"),\n    @(\"68\u00d796=6528\", \"12\u00d752=624\"),\n    @(\"61\u00d799=6039\", \"45\u00d741=1845\"),\n    @(\"46\u00d742=1932\", \"64\u00d748=3072\"),\n    @(\"68\u00d714=952\", \"22\u00d769=1518\"),\n    @(\"50\u00d727=1350\", \"76\u00d743=3268\"),\n    @(\"95\u00d784=7980\", \"25\u00d723=575\"),\n    @(\"69\u00d779=5451\", \"95\u00d716=1520\"),\n    @(\"86\u00d787=7482\", \"49\u00d722=1078\"),\n    @(\"92\u00d790=8280\", \"56\u00d740=2240\"),\n    @(\"74\u00d740=2960\", \"41\u00d796=3936\"),\n    @(\"35\u00d771=2485\", \"36\u00d717=612\"),\n    @(\"57\u00d742=2394\", \"30\u00d745=1350\"),\n    @(\"35\u00d722=770\", \"95\u00d720=1900\"),\n    @(\"38\u00d724=912\", \"63\u00d767=4221\"),\n    @(\"69\u00d790=6210\", \"17\u00d778=1326\"),\n    @(\"93\u00d740=3720\", \"44\u00d727=1188\"),\n    @(\"45\u00d768=3060\", \"35\u00d756=1960\"),\n    @(\"43\u00d745=1935\", \"58\u00d741=2378\"),\n    @(\"39\u00d797=3783\", \"61\u00d789=5429\"),\n    @(\"23\u00d724=552\", \"77\u00d790=6930\"),\n    @(\"49\u00d742=2058\", \"17\u00d784=1428\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
